# Auto-generated edit script applying numeric corrections to Sheets
# per the commit "chore: update Sheets via scheduled runner".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3667.7046
$ws.Range("I62").Value = 3052.2368
$ws.Range("K62").Value = 3052.2368
$ws.Range("M62").Value = -2428.2368

$ws.Range("H65").Value = 3667.7046
$ws.Range("I65").Value = 3052.2368
$ws.Range("K65").Value = 15261.184
$ws.Range("M65").Value = -12141.184

$ws.Range("H86").Value = 2003.5454
$ws.Range("I86").Value = 1927.4
$ws.Range("K86").Value = 1927.4
$ws.Range("M86").Value = -804.4000000000001

$ws.Range("H89").Value = 2003.5454
$ws.Range("I89").Value = 1927.4
$ws.Range("K89").Value = 9637
$ws.Range("M89").Value = -4021

$ws.Range("H103").Value = 1358.4667
$ws.Range("I103").Value = 2030.5
$ws.Range("J103").Value = 910.44446
$ws.Range("K103").Value = 6091.5
$ws.Range("L103").Value = 2731.33338
$ws.Range("M103").Value = -5505.5
$ws.Range("N103").Value = -3903.33338

$ws.Range("H112").Value = 5392.645
$ws.Range("J112").Value = 5786.5713
$ws.Range("L112").Value = 17359.7139
$ws.Range("N112").Value = -19575.7139

$ws.Range("H131").Value = 4100.1177
$ws.Range("J131").Value = 6299.8
$ws.Range("L131").Value = 18899.4
$ws.Range("N131").Value = -28979.4

$ws.Range("H132").Value = 5128887.5
$ws.Range("I132").Value = 5744215.5
$ws.Range("J132").Value = 1153.6666
$ws.Range("K132").Value = 17232646.5
$ws.Range("L132").Value = 3460.9998
$ws.Range("M132").Value = -17230116.5
$ws.Range("N132").Value = -8520.9998

$ws.Range("H137").Value = 12518.7
$ws.Range("J137").Value = 3543
$ws.Range("L137").Value = 10629
$ws.Range("N137").Value = -15729

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26605.488
$ws.Range("I32").Value = 28363.842
$ws.Range("K32").Value = 28363.842
$ws.Range("M32").Value = -28076.842

$ws.Range("H92").Value = 17505000
$ws.Range("J92").Value = 17505000
$ws.Range("L92").Value = 17505000
$ws.Range("N92").Value = -17509992

$ws.Range("H122").Value = 1375.0571
$ws.Range("J122").Value = 298
$ws.Range("L122").Value = 894
$ws.Range("N122").Value = -5794

$ws.Range("H132").Value = 1005.2353
$ws.Range("I132").Value = 851.3958
$ws.Range("J132").Value = 3466.6667
$ws.Range("K132").Value = 2554.1874
$ws.Range("L132").Value = 10400.0001
$ws.Range("M132").Value = -24.1873999999998
$ws.Range("N132").Value = -15460.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1388.2778
$ws.Range("I86").Value = 1350.4814
$ws.Range("K86").Value = 1350.4814
$ws.Range("M86").Value = -227.4813999999999

$ws.Range("H89").Value = 1388.2778
$ws.Range("I89").Value = 1350.4814
$ws.Range("K89").Value = 6752.406999999999
$ws.Range("M89").Value = -1136.406999999999

$ws.Range("H107").Value = 878.6
$ws.Range("I107").Value = 848.5
$ws.Range("K107").Value = 848.5
$ws.Range("M107").Value = 1071.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2440870
$ws.Range("I31").Value = 3031293.8
$ws.Range("J31").Value = 5372.125
$ws.Range("K31").Value = 3031293.8
$ws.Range("L31").Value = 5372.125
$ws.Range("M31").Value = -3030998.8
$ws.Range("N31").Value = -5962.125

$ws.Range("H34").Value = 2440870
$ws.Range("I34").Value = 3031293.8
$ws.Range("J34").Value = 5372.125
$ws.Range("K34").Value = 3031293.8
$ws.Range("L34").Value = 5372.125
$ws.Range("M34").Value = -3031091.8
$ws.Range("N34").Value = -5776.125

$ws.Range("H132").Value = 28795.611
$ws.Range("J132").Value = 1274.6
$ws.Range("L132").Value = 3823.8
$ws.Range("N132").Value = -8883.799999999999

$ws.Range("H134").Value = 2455.2
$ws.Range("I134").Value = 1940
$ws.Range("K134").Value = 5820
$ws.Range("M134").Value = -3285

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 333.1
$ws.Range("J38").Value = 464.42856
$ws.Range("L38").Value = 1393.28568
$ws.Range("N38").Value = -2087.28568

$ws.Range("H92").Value = 1857.7
$ws.Range("I92").Value = 1822.125
$ws.Range("J92").Value = 2000
$ws.Range("K92").Value = 5466.375
$ws.Range("L92").Value = 6000
$ws.Range("M92").Value = -4218.375
$ws.Range("N92").Value = -8496

$ws.Range("H129").Value = 1930.3914
$ws.Range("I129").Value = 1537.25
$ws.Range("J129").Value = 2359.2727
$ws.Range("K129").Value = 4611.75
$ws.Range("L129").Value = 7077.8181
$ws.Range("M129").Value = 388.25
$ws.Range("N129").Value = -17077.8181

$ws.Range("H140").Value = 3818.8462
$ws.Range("I140").Value = 3818.8462
$ws.Range("K140").Value = 11456.5386
$ws.Range("M140").Value = -6276.5386

$ws.Range("H141").Value = 5415.6
$ws.Range("I141").Value = 4154.636
$ws.Range("K141").Value = 12463.908
$ws.Range("M141").Value = -7283.908000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 13227
$ws.Range("I80").Value = 1970
$ws.Range("J80").Value = 19981.2
$ws.Range("K80").Value = 1970
$ws.Range("L80").Value = 19981.2
$ws.Range("M80").Value = -972
$ws.Range("N80").Value = -21977.2

$ws.Range("H83").Value = 13227
$ws.Range("I83").Value = 1970
$ws.Range("J83").Value = 19981.2
$ws.Range("K83").Value = 9850
$ws.Range("L83").Value = 99906
$ws.Range("M83").Value = -4858
$ws.Range("N83").Value = -109890

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3382.5
$ws.Range("I7").Value = 3618.6
$ws.Range("K7").Value = 3618.6
$ws.Range("M7").Value = -3506.6

$ws.Range("H55").Value = 1010.93335
$ws.Range("I55").Value = 247
$ws.Range("K55").Value = 247
$ws.Range("M55").Value = -74

$ws.Range("H82").Value = 1818.8
$ws.Range("I82").Value = 1527.1428
$ws.Range("K82").Value = 1527.1428
$ws.Range("M82").Value = -1166.1428

$ws.Range("H85").Value = 1818.8
$ws.Range("I85").Value = 1527.1428
$ws.Range("K85").Value = 1527.1428
$ws.Range("M85").Value = -279.1428000000001

$ws.Range("H126").Value = 3382.5
$ws.Range("I126").Value = 3618.6
$ws.Range("K126").Value = 10855.8
$ws.Range("M126").Value = -8385.799999999999

$ws.Range("H132").Value = 2916.3333
$ws.Range("I132").Value = 1997
$ws.Range("K132").Value = 5991
$ws.Range("M132").Value = -3461

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 7761
$ws.Range("I23").Value = 1500
$ws.Range("K23").Value = 1500
$ws.Range("M23").Value = -1271

$ws.Range("H100").Value = 1432.6
$ws.Range("I100").Value = 1322.25
$ws.Range("K100").Value = 2644.5
$ws.Range("M100").Value = -2103.5

$ws.Range("H107").Value = 858.2727
$ws.Range("J107").Value = 757.6
$ws.Range("L107").Value = 2272.8
$ws.Range("N107").Value = -6112.8

$ws.Range("H126").Value = 2777.5833
$ws.Range("I126").Value = 1944.1428
$ws.Range("J126").Value = 3944.4
$ws.Range("K126").Value = 5832.428400000001
$ws.Range("L126").Value = 11833.2
$ws.Range("M126").Value = -3362.428400000001
$ws.Range("N126").Value = -16773.2
